# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de),
# representing the "ded75eeb-...md" file that is now "Ready for handoff" / has a new
# handoff xliff generated, wires up the matching hyperlinks, resizes the tables, and
# widens a couple of columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # Overview
$ws2 = $wb.Worksheets.Item(2)  # zh-cn
$ws3 = $wb.Worksheets.Item(3)  # de-de

$dateFormat = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276   # BGR encoding of RGB(0x64,0x95,0xED) -> matches existing HyperLink style color FF6495ED

# ============================================================
# Sheet "Overview": add row 3
# ============================================================
$ws1.Range("A3").Value = 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ws1.Range("B3").Value = 'e2e\ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ws1.Range("C3").Value = '.md'
$ws1.Range("D3").Value = ''
$ws1.Range("E3").Value = 'Ready for handoff'
$ws1.Range("F3").Value = 'Ready for handoff'
$ws1.Range("G3").Value = '2017-02-09 10:38:38'

# Style B3 like the existing hyperlink cell (B2) and G3 like the existing date cell (G2)
$ws1.Range("B3").Font.Underline = 2
$ws1.Range("B3").Font.Color = $hyperlinkColor
$ws1.Range("G3").NumberFormat = $dateFormat

# Hyperlink for B3 (mirrors the hyperlink already present on B2)
$ws1.Hyperlinks.Add($ws1.Range("B3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bf4ce38c8ca9c8f7d01473e0cfe372d431463bd/e2e/ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', "", "", 'e2e\ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md') | Out-Null

# Widen columns E (zh-cn) and F (de-de) to fit the new status text
$ws1.Columns.Item(5).ColumnWidth = 16.3
$ws1.Columns.Item(6).ColumnWidth = 16.3

# Grow the Overview table to include the new row
$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G3"))

# ============================================================
# Sheet "zh-cn": add row 3
# ============================================================
$ws2.Range("A3").Value = 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ws2.Range("B3").Value = '.md'
$ws2.Range("C3").Value = 'Ready for handoff'
$ws2.Range("D3").Value = 'e2e'
$ws2.Range("E3").Value = 'ht'
$ws2.Range("F3").Value = 'False'
$ws2.Range("G3").Value = 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5fooooooooooooooooooooooooooooooooooooooo.190c1e32e02e017e00acdf087b57c98c68aee04c.zh-cn.xlf'
$ws2.Range("H3").Value = '2017-02-09 10:38:16'
$ws2.Range("I3").Value = ''
$ws2.Range("J3").Value = ''
$ws2.Range("K3").Value = ''
$ws2.Range("L3").Value = '0001-01-01 00:00:00'
$ws2.Range("M3").Value = ''
$ws2.Range("N3").Value = ''
$ws2.Range("O3").Value = 'True'
$ws2.Range("P3").Value = ''
$ws2.Range("Q3").Value = 'False'
$ws2.Range("R3").Value = ''

# Style A3 like the existing hyperlink cell (A2) and H3/L3 like the existing date cells (H2/L2)
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = $hyperlinkColor
$ws2.Range("H3").NumberFormat = $dateFormat
$ws2.Range("L3").NumberFormat = $dateFormat

# Hyperlink for A3 (mirrors the hyperlink already present on A2)
$ws2.Hyperlinks.Add($ws2.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bf4ce38c8ca9c8f7d01473e0cfe372d431463bd/e2e/ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', "", "", 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md') | Out-Null

# Widen column C to fit the new status text
$ws2.Columns.Item(3).ColumnWidth = 16.3

# Grow the zh-cn table to include the new row
$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:R3"))

# ============================================================
# Sheet "de-de": add row 3
# ============================================================
$ws3.Range("A3").Value = 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ws3.Range("B3").Value = '.md'
$ws3.Range("C3").Value = 'Ready for handoff'
$ws3.Range("D3").Value = 'e2e'
$ws3.Range("E3").Value = 'ht'
$ws3.Range("F3").Value = 'False'
$ws3.Range("G3").Value = 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5fooooooooooooooooooooooooooooooooooooooo.190c1e32e02e017e00acdf087b57c98c68aee04c.de-de.xlf'
$ws3.Range("H3").Value = '2017-02-09 10:38:38'
$ws3.Range("I3").Value = ''
$ws3.Range("J3").Value = ''
$ws3.Range("K3").Value = ''
$ws3.Range("L3").Value = '0001-01-01 00:00:00'
$ws3.Range("M3").Value = ''
$ws3.Range("N3").Value = ''
$ws3.Range("O3").Value = 'True'
$ws3.Range("P3").Value = ''
$ws3.Range("Q3").Value = 'False'
$ws3.Range("R3").Value = ''

# Style A3 like the existing hyperlink cell (A2) and H3/L3 like the existing date cells (H2/L2)
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = $hyperlinkColor
$ws3.Range("H3").NumberFormat = $dateFormat
$ws3.Range("L3").NumberFormat = $dateFormat

# Hyperlink for A3 (mirrors the hyperlink already present on A2)
$ws3.Hyperlinks.Add($ws3.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bf4ce38c8ca9c8f7d01473e0cfe372d431463bd/e2e/ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', "", "", 'ded75eeb-ebe0-4c31-a25d-7cb06af6ca5foooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md') | Out-Null

# Widen column C to fit the new status text
$ws3.Columns.Item(3).ColumnWidth = 16.3

# Grow the de-de table to include the new row
$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:R3"))

